# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型"
# sheets to reflect the latest generated data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Map of known old values -> new values for column F ("想去人数")
    $updates = @{
        730  = 731
        2747 = 2750
        3717 = 3720
    }

    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        $cell = $ws.Cells.Item($r, 6)  # Column F
        $val = $cell.Value2
        if ($null -ne $val -and $updates.ContainsKey([int]$val)) {
            $cell.Value = $updates[[int]$val]
        }
    }
}

$wb.Save()
